$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 4, pushing the existing rows 4-29 down to 6-31
$ws.Rows("4:5").Insert()

# Apply the style used by column A (index style) to the new rows' A cells
$ws.Range("A4").Value = 2
$ws.Range("A4").Style = $ws.Range("A6").Style
$ws.Range("B4").Value = "Holden"

$ws.Range("A5").Value = 3
$ws.Range("A5").Style = $ws.Range("A6").Style
$ws.Range("B5").Value = "Rizzie Spiral"

# New data values for row 4 ("Holden")
$row4vals = @(0.9455870434453302, 0.9259054374446701, 1.11752398176255, 0.8809408136266945, 1.224048710175232, 0.9259054374446701, 0.9537535868607189, 1.224048710175232, 0.9259054374446701, 0.9537535868607189, 1.088901148517976, 1.088901148517976, 1.098442092932834, 1.034569244826874, 1.034569244826874, 1.007403292981323, 1.007403292981323, 1.007959928885866)

# New data values for row 5 ("Rizzie Spiral")
$row5vals = @(1.104084103176118, 1.119535777142581, 0.7577003905768592, 1.317950710269136, 0.5926340241291361, 1.119535777142581, 1.096097942123474, 0.5926340241291361, 1.119535777142581, 1.096097942123474, 0.844365983126305, 0.844365983126305, 0.8154774522764897, 0.9360892477983969, 0.9360892477983969, 0.9819508801344429, 0.9819508801344429, 0.9980004912362174)

for ($i = 0; $i -lt $row4vals.Length; $i++) {
    $ws.Cells.Item(4, 3 + $i).Value = $row4vals[$i]
}

for ($i = 0; $i -lt $row5vals.Length; $i++) {
    $ws.Cells.Item(5, 3 + $i).Value = $row5vals[$i]
}

# Rename the shared string "Thomas Hex" -> "Matthies Hex" (now located in what was row 9, now row 11)
$ws.Range("B11").Value = "Matthies Hex"
